$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 headers for the two new iteration blocks ---
$ws.Range("E1").Value = "Iteration_1"
$ws.Range("H1").Value = "Iteration_2"
$ws.Range("E1:G1").Merge()
$ws.Range("H1:J1").Merge()

# --- Row 2 year labels for the new blocks ---
$ws.Range("E2").Value = "2030"
$ws.Range("F2").Value = "2040"
$ws.Range("G2").Value = "2050"
$ws.Range("H2").Value = "2030"
$ws.Range("I2").Value = "2040"
$ws.Range("J2").Value = "2050"

# --- Update existing Standalone values (rows 4-6) ---
$ws.Range("B4").Value = 1184000.000000001
$ws.Range("C4").Value = 4810.151102347427
$ws.Range("D4").Value = [double]"-2.620601923793058e-10"

$ws.Range("C5").Value = 1175482.441022393
$ws.Range("D5").Value = 1180214.184748082

$ws.Range("D6").Value = 71.29135345365613

# --- Iteration_1 (E:G) and Iteration_2 (H:J) values ---
$ws.Range("E4").Value = 1184000.000000074
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 1183999.999999876
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0

$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0

$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 1181739.817836299
$ws.Range("G6").Value = 1181730.993934782
$ws.Range("H6").Value = [double]"4.984140500134267e-07"
$ws.Range("I6").Value = 1181737.653169923
$ws.Range("J6").Value = 1181730.993934819

$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = [double]"6.694979318452933e-09"
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0

# Rows 8-15: all-zero fill for both new iteration blocks
for ($r = 8; $r -le 15; $r++) {
    $ws.Range("E$r").Value = 0
    $ws.Range("F$r").Value = 0
    $ws.Range("G$r").Value = 0
    $ws.Range("H$r").Value = 0
    $ws.Range("I$r").Value = 0
    $ws.Range("J$r").Value = 0
}
